$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2025-02-02 Sunday" "2025-02-03 Monday"

Replace-Text "400÷2=200, 0" "960÷7=137, 1"
Replace-Text "947÷7=135, 2" "911÷2=455, 1"
Replace-Text "859÷2=429, 1" "532÷3=177, 1"
Replace-Text "682÷9=75, 7" "880÷3=293, 1"
Replace-Text "305÷5=61, 0" "904÷9=100, 4"

Replace-Text "804÷8=100, 4" "843÷2=421, 1"
Replace-Text "418÷4=104, 2" "679÷5=135, 4"
Replace-Text "364÷7=52, 0" "264÷3=88, 0"
Replace-Text "285÷3=95, 0" "135÷4=33, 3"
Replace-Text "899÷9=99, 8" "250÷6=41, 4"

Replace-Text "224÷8=28, 0" "899÷5=179, 4"
Replace-Text "545÷9=60, 5" "343÷6=57, 1"
Replace-Text "732÷6=122, 0" "315÷8=39, 3"
Replace-Text "707÷8=88, 3" "909÷2=454, 1"
Replace-Text "681÷6=113, 3" "544÷9=60, 4"

Replace-Text "503÷4=125, 3" "643÷9=71, 4"
Replace-Text "420÷5=84, 0" "121÷4=30, 1"
Replace-Text "937÷5=187, 2" "861÷2=430, 1"
Replace-Text "549÷7=78, 3" "638÷7=91, 1"
Replace-Text "716÷5=143, 1" "516÷5=103, 1"

Replace-Text "161÷2=80, 1" "513÷3=171, 0"
Replace-Text "195÷7=27, 6" "411÷9=45, 6"
Replace-Text "992÷6=165, 2" "349÷5=69, 4"
Replace-Text "979÷6=163, 1" "253÷9=28, 1"
Replace-Text "281÷8=35, 1" "320÷7=45, 5"
